$d = $word.ActiveDocument

# --- Step 1: Replace "Estrategia de Pesquisa" paragraph text ---
$oldText = "Descrever a estratégia e etiquetagem (labeling) utilizada ou implementada, nomeadamente no que diz respeito à ordenação de variáveis e valores."
$newText = "A estratégia de pesquisa consiste nas restrições serem aplicadas da seguinte forma: primeiro não permitindo que existam espaços 2x2 com peças da mesma core e depois verificando se todas as peças estão ligadas a pelo menos outra da sua cor."

$found = $false
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "*Descrever a estratégia*") {
        $r.Text = $newText
        $r.Font.Size = 11
        $found = $true
        break
    }
}
Write-Output "step1 found: $found"
